# Update TPM-derived values in the active worksheet (Mmp12-Plaur LR-pairs data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.007057333333333333
$ws.Range("H2").Value = 0.021172
$ws.Range("M2").Value = 2.843949
$ws.Range("N2").Value = 8.531846999999999
$ws.Range("O2").Value = 0.4976240243095911
$ws.Range("P2").Value = 0.4976240243095912
$ws.Range("Q2").Value = 0.020070696076
$ws.Range("R2").Value = 0.180636264684
$ws.Range("S2").Value = 0.4976240243095911
$ws.Range("T2").Value = 0.4976240243095912

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.007057333333333333
$ws.Range("H3").Value = 0.021172
$ws.Range("O3").Value = 0.4403664892852895
$ws.Range("P3").Value = 0.4403664892852897
$ws.Range("Q3").Value = 0.01776132488933333
$ws.Range("R3").Value = 0.159851924004
$ws.Range("S3").Value = 0.4403664892852895
$ws.Range("T3").Value = 0.4403664892852897

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.007057333333333333
$ws.Range("H4").Value = 0.021172
$ws.Range("M4").Value = 0.3543876666666666
$ws.Range("O4").Value = 0.06200948640511928
$ws.Range("P4").Value = 0.0620094864051193
$ws.Range("Q4").Value = 0.002501031892888888
$ws.Range("R4").Value = 0.022509287036
$ws.Range("S4").Value = 0.06200948640511928
$ws.Range("T4").Value = 0.0620094864051193
